$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Jimmy Butler III", "SF,PF", "Golden State Warriors"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Jonas Valanciunas", "C", "Sacramento Kings"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
